# ----------------------------------------------------------------------------
# Rebuild "Solicitudes" sheet: expand the 8-column sample-request template
# into a 21-column template (Age ... Monthly_Balance) with 3 sample rows,
# refresh the instructions banner text, and adjust the view a little.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Instructions banner (row 1, currently merged A1:H1) -----------------
$ws.Range("A1").Value = "Instrucciones: Reemplace los valores de ejemplo con los datos del nuevo cliente y guarde el archivo antes de ejecutar la FASE 11."

# --- 2. Header row (row 2) ---------------------------------------------------
$ws.Range("A2").Value  = "Age"
$ws.Range("B2").Value  = "Occupation"
$ws.Range("C2").Value  = "Annual_Income"
$ws.Range("D2").Value  = "Monthly_Inhand_Salary"
$ws.Range("E2").Value  = "Num_Bank_Accounts"
$ws.Range("F2").Value  = "Num_Credit_Card"
$ws.Range("G2").Value  = "Interest_Rate"
$ws.Range("H2").Value  = "Num_of_Loan"
$ws.Range("I2").Value  = "Delay_from_due_date"
$ws.Range("J2").Value  = "Num_of_Delayed_Payment"
$ws.Range("K2").Value  = "Changed_Credit_Limit"
$ws.Range("L2").Value  = "Num_Credit_Inquiries"
$ws.Range("M2").Value  = "Credit_Mix"
$ws.Range("N2").Value  = "Outstanding_Debt"
$ws.Range("O2").Value  = "Credit_Utilization_Ratio"
$ws.Range("P2").Value  = "Credit_History_Age"
$ws.Range("Q2").Value  = "Payment_of_Min_Amount"
$ws.Range("R2").Value  = "Total_EMI_per_month"
$ws.Range("S2").Value  = "Amount_invested_monthly"
$ws.Range("T2").Value  = "Payment_Behaviour"
$ws.Range("U2").Value  = "Monthly_Balance"

# New header cells (I2:U2) need the same look as the rest of row 2 (style
# used by A2:H2). Copy that formatting across without disturbing the values
# we just wrote, and without introducing a brand-new style record.
$ws.Range("A2").Copy()
$ws.Range("I2:U2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- 3. Sample data rows 3-5 -------------------------------------------------
# Row 3
$ws.Cells.Item(3, 1).Value  = 35
$ws.Cells.Item(3, 3).Value  = 480000
$ws.Cells.Item(3, 4).Value  = 40000
$ws.Cells.Item(3, 5).Value  = 3
$ws.Cells.Item(3, 6).Value  = 2
$ws.Cells.Item(3, 7).Value  = 12
$ws.Cells.Item(3, 8).Value  = 1
$ws.Cells.Item(3, 9).Value  = 2
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 5
$ws.Cells.Item(3, 12).Value = 3
$ws.Cells.Item(3, 15).Value = 28
$ws.Cells.Item(3, 16).Value = 100
$ws.Cells.Item(3, 18).Value = 2000
$ws.Cells.Item(3, 19).Value = 1500
$ws.Cells.Item(3, 21).Value = 5000

# Row 4
$ws.Cells.Item(4, 1).Value  = 42
$ws.Cells.Item(4, 3).Value  = 1080000
$ws.Cells.Item(4, 4).Value  = 85000
$ws.Cells.Item(4, 5).Value  = 4
$ws.Cells.Item(4, 6).Value  = 2
$ws.Cells.Item(4, 7).Value  = 9
$ws.Cells.Item(4, 8).Value  = 1
$ws.Cells.Item(4, 9).Value  = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 15).Value = 15
$ws.Cells.Item(4, 16).Value = 220
$ws.Cells.Item(4, 18).Value = 3000
$ws.Cells.Item(4, 19).Value = 10000
$ws.Cells.Item(4, 21).Value = 15000

# Row 5
$ws.Cells.Item(5, 1).Value  = 23
$ws.Cells.Item(5, 3).Value  = 180000
$ws.Cells.Item(5, 4).Value  = 12000
$ws.Cells.Item(5, 5).Value  = 1
$ws.Cells.Item(5, 6).Value  = 3
$ws.Cells.Item(5, 7).Value  = 22
$ws.Cells.Item(5, 8).Value  = 4
$ws.Cells.Item(5, 9).Value  = 10
$ws.Cells.Item(5, 10).Value = 7
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 9
$ws.Cells.Item(5, 14).Value = 65000
$ws.Cells.Item(5, 15).Value = 85
$ws.Cells.Item(5, 16).Value = 70
$ws.Cells.Item(5, 18).Value = 7000
$ws.Cells.Item(5, 19).Value = 200
$ws.Cells.Item(5, 21).Value = 500

# --- 4. Column widths (extend the existing 20-wide columns out to U) --------
# ColumnWidth is in characters; 19.2 chars round-trips to the same stored
# <col width="20"/> already used by columns A:H.
$ws.Range("I1:U1").ColumnWidth = 19.2

# --- 5. Re-merge the banner across the new column range ---------------------
$ws.Range("A1:U1").Merge()
# Merge() re-stamps every cell in the range with the anchor's style; put the
# non-anchor cells (B1:U1) back to their original (unformatted) look.
$ws.Range("B1:U1").ClearFormats()

# --- 6. View tweaks (zoom + selection) --------------------------------------
$ws.Application.ActiveWindow.Zoom = 186
$ws.Range("N3").Select()
